# Correct the Activity / PlannedResource values for the "Front Cover   2p"
# rows (4-5) and the "Text  2p" rows (6-7): each pair had its Digital Print
# and Proof Approval steps swapped; put them back in the right rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Front Cover   2p  (rows 4-5)
$ws.Range("B4").Value = "Proof Approval"
$ws.Range("D4").Value = "STL Proof Approval"
$ws.Range("L4").Value = "STL Proof Approval"
$ws.Range("M4").Value = "STL Proof Approval"

$ws.Range("B5").Value = "Digital Print F 4x0"
$ws.Range("D5").Value = "HC NexPress 1-4c"
$ws.Range("L5").Value = "HC NexPress 1-4c"
$ws.Range("M5").Value = "HC NexPress 1-4c"

# Text  2p  (rows 6-7)
$ws.Range("B6").Value = "Proof Approval"
$ws.Range("D6").Value = "STL Proof Approval"
$ws.Range("L6").Value = "STL Proof Approval"
$ws.Range("M6").Value = "STL Proof Approval"

$ws.Range("B7").Value = "Digital Print F/B 4x4"
$ws.Range("D7").Value = "HC NexPress 1-4c"
$ws.Range("L7").Value = "HC NexPress 1-4c"
$ws.Range("M7").Value = "HC NexPress 1-4c"
